# Generate Report for Archive
# - Update localization status text from "Ready for handoff" to "In Translation"
#   on the Overview sheet (zh-cn / de-de status columns) as well as the
#   per-locale detail sheets (zh-cn, de-de).
# - Narrow the now-shorter "Status" columns to their new (re-fitted) width.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Overview sheet: columns E (zh-cn) and F (de-de) hold the status per file.
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# Per-locale detail sheets: column C holds the Status value.
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# The shorter status text means the Status columns can shrink; match the
# narrower re-fitted width that results (stored column width ends up at
# 13.4101845877511 characters).
$newColumnWidth = 13.4101845877511 - (5 / 6)

$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth
